$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.42809999999999
$ws.Range("C3").Value = -11.2174
$ws.Range("D5").Value = -8.370199999999997
$ws.Range("C14").Value = -12.6575
$ws.Range("C16").Value = -12.62410000000001
$ws.Range("D16").Value = -8.154300000000008
$ws.Range("C21").Value = -13.24170000000001
$ws.Range("C23").Value = -12.0882
$ws.Range("C25").Value = -11.2483
